$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.398.54'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.844.73'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6313'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07538'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2927'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07711'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '1.837.41'
$ws.Range("E12").Value = '  -7.44%  '
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6800'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("E15").Value = '  +5.35%  '
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '2.087.52'
$ws.Range("E17").Value = '  -7.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.166'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '29.431.80'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.453'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.355'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.460'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.287'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.022'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.848'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7104'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.592'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").Value = '1.246.85'
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01808'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.384'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9019'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000118'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.093'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.3998'
$ws.Range("D48").Style = "Normal"
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.673'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.923'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1121'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.54%  '
